# Add additional regional data for steel production to the "connections" sheet.
# A new row is inserted after the existing "steel <- pellets (simple_BF)" row,
# describing a "primary fuel" connection, and the two rows that follow it
# (previously the "coke" rows for steel/simple_BF-all and pellets/all) are
# updated to describe a "fossil fuel" connection produced via simple_sinter /
# simple_pellets instead of the generic "all" unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Insert a new row at row 3 - this shifts the old rows 3.. down to 4.. and
# copies the formatting of row 2 (the row above) into the fresh row, just
# like Excel's normal "Insert" behaviour.
$ws.Rows(3).Insert()

# --- Row 4 (was row 3 before the insert): switch its unit + inflow product
# to describe "fossil fuel" sourced through simple_sinter.
$ws.Range("C4").Value = "simple_sinter"
$ws.Range("E4").Value = "fossil fuel"
$ws.Range("G4").Value = "outflows"

# --- Row 5 (was row 4 before the insert): switch its unit + inflow product
# to describe "fossil fuel" sourced through simple_pellets, and pick up the
# "Text" number format (no bold font flag) used by the other unit cells.
$ws.Range("C5").Value = "simple_pellets"
$ws.Range("H5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("E5").Value = "fossil fuel"
$ws.Range("F5").ClearFormats()

# --- New row 3: the "primary fuel" connection for steel via simple_BF.
$ws.Range("B3").Value = "steel"
$ws.Range("C3").Value = "simple_BF"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "inflow"
$ws.Range("E3").Value = "primary fuel"
$ws.Range("F3").Value = "coke"
$ws.Range("G3").Value = "outflow"
$ws.Range("H3").Value = "simple_coke"
$ws.Range("I3").Value = "coke"

# Leave the cursor on the newly-entered cell, matching the saved selection.
$ws.Range("E3").Select()
